# [Outlook] (mapping) Include new snippets
# Append new rows to the "Snippets" table describing the additional
# "Other item APIs" mappings (subject, internetMessageId, itemClass,
# itemType, start) for Outlook add-in snippets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$rows = @(
    @("AppointmentRead",    "subject",            $null, "outlook-other-item-apis-get-subject-read",                    "get"),
    @("MessageRead",        "subject",            $null, "outlook-other-item-apis-get-subject-read",                    "get"),
    @("AppointmentCompose", "subject",            $null, "outlook-other-item-apis-get-set-subject-compose",             "get"),
    @("MessageCompose",     "subject",            $null, "outlook-other-item-apis-get-set-subject-compose",             "get"),
    @("AppointmentCompose", "subject",            $null, "outlook-other-item-apis-get-set-subject-compose",             "set"),
    @("MessageCompose",     "subject",            $null, "outlook-other-item-apis-get-set-subject-compose",             "set"),
    @("MessageRead",        "internetMessageId",  $null, "outlook-other-item-apis-get-internet-message-id-read",        "get"),
    @("AppointmentRead",    "itemClass",          $null, "outlook-other-item-apis-get-item-class-read",                 "get"),
    @("MessageRead",        "itemClass",          $null, "outlook-other-item-apis-get-item-class-read",                 "get"),
    @("AppointmentCompose", "itemType",           $null, "outlook-other-item-apis-get-item-type",                       "get"),
    @("AppointmentRead",    "itemType",           $null, "outlook-other-item-apis-get-item-type",                       "get"),
    @("MessageCompose",     "itemType",           $null, "outlook-other-item-apis-get-item-type",                       "get"),
    @("MessageRead",        "itemType",           $null, "outlook-other-item-apis-get-item-type",                       "get"),
    @("AppointmentRead",    "start",              $null, "outlook-other-item-apis-get-start-read",                      "get"),
    @("MessageRead",        "start",              $null, "outlook-other-item-apis-get-start-read",                      "get"),
    @("AppointmentCompose", "start",              $null, "outlook-other-item-apis-get-set-start-appointment-organizer", "get"),
    @("Time",               "getAsync",           2,     "outlook-other-item-apis-get-set-start-appointment-organizer", "get"),
    @("AppointmentCompose", "start",              $null, "outlook-other-item-apis-get-set-start-appointment-organizer", "set"),
    @("Time",               "setAsync",           2,     "outlook-other-item-apis-get-set-start-appointment-organizer", "set")
)

foreach ($row in $rows) {
    $newRow = $lo.ListRows.Add()
    $r = $lo.Range.Row + $lo.Range.Rows.Count - 1

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("A2").Select() | Out-Null
